# Penalty Reward System (unfinished) - shift forecast weeks forward by one
# week and refresh the dependent Summary statistics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Forecast Comparison": every Week_Start_Date (col B) moves one
# week later. Pre-format column B as Text so Excel does not silently
# reinterpret the literal "YYYY-MM-DD" strings as date serial numbers.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("B2:B17").NumberFormat = "@"

$ws1.Range("B2").Value = "2025-01-12"
$ws1.Range("B3").Value = "2025-01-19"
$ws1.Range("B4").Value = "2025-01-26"
$ws1.Range("B5").Value = "2025-02-02"
$ws1.Range("B6").Value = "2025-02-09"
$ws1.Range("B7").Value = "2025-02-16"
$ws1.Range("B8").Value = "2025-02-23"
$ws1.Range("B9").Value = "2025-03-02"
$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("B17").Value = "2025-04-27"

# MyForecast values that changed along with the week shift.
$ws1.Range("D8").Value = 5
$ws1.Range("D9").Value = 6

# ---------------------------------------------------------------------
# Sheet "Summary": figures recomputed from the updated forecast above.
# Pre-format column B as Text so numeric-looking strings ("73", "37", ...)
# and the date-looking strings stay literal text, matching the source.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B8:B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2024-02-11 to 2025-01-05"
$ws2.Range("B8").Value = "211 units"
$ws2.Range("B9").Value = "73"
$ws2.Range("B10").Value = "37"
$ws2.Range("B11").Value = "18"
$ws2.Range("B12").Value = "6"
$ws2.Range("B13").Value = "2025-03-02"
$ws2.Range("B14").Value = "4"
$ws2.Range("B15").Value = "2025-01-12"
